$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new backlog items right after the existing "Collapse/expand" row,
# before touching the rows whose shared-string content changes. This keeps the
# shared string table growth in the same order as the source edit.
$ws.Range("A8").Value = "alternate project/task delete. in project properties panel"
$ws.Range("A9").Value = "better tree navigation system: hand grab and drag, diagonal arrow, more intuitive"
$ws.Range("A10").Value = "Project progress incorporated visually into tree/outline views"

# Fix the typo ("on;y" -> "only") in the project hours item and move it down one row.
$ws.Range("A6").Value = "Project hours invested. So far the project is defined only in calendar time. There will also be time spent/invested in the project"

# The "update hours spent" item moves up to row 7 (same text as before).
$ws.Range("A7").Value = "once the project hours feature is implemented, the user will need to be able to update the hours spent."

# Append the remaining new backlog items.
$ws.Range("A11").Value = "when project/task is completed, fireworks, music, congratulations"
$ws.Range("A12").Value = "project accomplishments share on facebook"
$ws.Range("A13").Value = "project documentation - add pictures of project in progress and on completion"
$ws.Range("A14").Value = "project archiving"

$ws.Range("A14").Select()
